$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New loss-of-sale records to append starting at row 19
$rows = @(
    @{ A=17; B="22-12-2025"; C="FAVAS";  D=7356899762; E="24-12-2025"; F="MUHAMMED THUFAIL C I"; G="Loss"; H="PRICING"; I="RENT TO HIGH"; J="-"; K="PRICING" },
    @{ A=18; B="22-12-2025"; C="AFNAS";  D=9778549640; E="01-01-2026"; F="Abdul Hadi Rafeeque"; G="Loss"; H="ENQUIRY"; I="Enquiry for Relative/Friend"; J="-"; K="ENQUIRY" },
    @{ A=19; B="23-12-2025"; C="Roshan"; D=9544282330; E="23-12-2025"; F="Abdul Hadi Rafeeque"; G="Loss"; H="PRODUCT"; I="PRODUCT NOT AVAILABLE"; J="-"; K="PRODUCT ALREADY BOOKED AS A SAME DATE ANOTHER CUSTOMER" },
    @{ A=20; B="24-12-2025"; C="ARSHAD";  D=9946858090; E="07-02-2026"; F="SHIBIN RAJ KK"; G="Loss"; H="SIZE NOT SUITABLE"; I="SIZE TOO SMALL"; J="-"; K=$null },
    @{ A=21; B="24-12-2025"; C="RIYAS";   D=9947154435; E="04-01-2026"; F="MUHAMMED THUFAIL C I"; G="Loss"; H="ENQUIRY"; I="-"; J="-"; K="JUST ENQUIRY" },
    @{ A=22; B="25-12-2025"; C="Boby";    D=9847503811; E="27-12-2025"; F="SHIBIN RAJ KK"; G="Loss"; H="ENQUIRY"; I="Enquiry for Relative/Friend"; J="-"; K="tomorrow will come" },
    @{ A=23; B="25-12-2025"; C="IRSHAD";  D=9645503686; E="10-01-2026"; F="Abdul Hadi Rafeeque"; G="Loss"; H="PRODUCT"; I="REQUIRED DESIGN NOT AVAILABLE"; J="-"; K="ENQUIRY" }
)

$startRow = 19
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data.A
    $ws.Cells.Item($r, 1).Style = $ws.Cells.Item($r - 1, 1).Style

    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C

    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 4).Style = $ws.Cells.Item($r - 1, 4).Style

    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $data.F
    $ws.Cells.Item($r, 7).Value = $data.G
    $ws.Cells.Item($r, 8).Value = $data.H
    $ws.Cells.Item($r, 9).Value = $data.I
    $ws.Cells.Item($r, 10).Value = $data.J

    if ($null -ne $data.K) {
        $ws.Cells.Item($r, 11).Value = $data.K
    }
}
